$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final BOM rows (2..17), pipe-delimited:
# Description|DesignItemId|Designator|Footprint|LibRef|Quantity
$rows = @(
  "Capacitor|C0805C104K5RAC7411|C1, C2, C3, C6, C7, C10, C11|CAPC2012X94N|C0805C104K5RAC7411|7",
  "Capacitor|08055A4R3CAT2A|C4, C5|CAPC2012X94N|08055A4R3CAT2A|2",
  "Capacitor|885012107015|C8, C9|CAPC2012X135N|885012107015|2",
  "USB connector|920-E52A2021S10100|CN1|USB307530A|920-E52A2021S10100|1",
  "Integrated Circuit|STM32L031K6T6|IC1|QFP80P900X900X160-32N|STM32L031K6T6|1",
  "Power Supply|BU33SD5WG-TR|IC2|SOT95P280X125-5N|BU33SD5WG-TR|1",
  "Integrated Circuit|TB67H450FNG,EL|IC3|SOIC127P600X175-9N|TB67H450FNG,EL|1",
  "Connector|4PIN|J1|4PIN|4PIN|1",
  "Connector|B2B-PH-K-S_LF__SN_|J2|SHDR2W50P0X200_1X2_590X450X600P|B2B-PH-K-S_LF__SN_|1",
  "Inductor|742792030|L1|RESC2012X60N|742792030|1",
  "LED|APT2012LZGCK|LED1|LEDC2012X85N|APT2012LZGCK|1",
  "Resistor|RC0805FR-071K8L|R1|RESC2012X60N|RC0805FR-071K8L|1",
  "Resistor|RC0805FR-0710KL|R2|RESC2012X60N|RC0805FR-0710KL|1",
  "Jumper (0 ohm)|RC0805FR-070RL|SB1, SB2, SB3, SB4, SB5|RESC2012X60N|RC0805FR-070RL|5",
  "Tactile switch|EVP-BT3G4A000|SW1|EVP-BT3G4A000|EVP-BT3G4A000|1",
  "Crystal|ABS07-32.768KHZ-9-T|X1|ABS-07|ABS07-32.768KHZ-9-T|1"
)

# Stash the original text-cell format (style "2": thin border + quote-prefix)
# and number-cell format (style "1": thin border) in scratch cells before we
# start overwriting A2:F17 -- writing .Value into those source cells later
# would otherwise destroy the very formatting we want to copy from them.
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(2, 6).Copy() | Out-Null
$ws.Range("H2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$startRow = 2
for ($i = 0; $i -lt $rows.Length; $i++) {
    $parts = $rows[$i].Split("|")
    $rnum = $startRow + $i
    for ($c = 1; $c -le 5; $c++) {
        $val = $parts[$c - 1]
        if ($val -match '^[0-9]+$') {
            # Purely-numeric DesignItemId/LibRef text (e.g. "885012107015")
            # must stay text, matching the source data, instead of being
            # auto-coerced to a number -- force it with a leading quote.
            $ws.Cells.Item($rnum, $c).Value = "'" + $val
        } else {
            $ws.Cells.Item($rnum, $c).Value = $val
        }
    }
    $qty = [int]$parts[5]
    $ws.Cells.Item($rnum, 6).Value = $qty
}

# Re-apply the stashed formatting (border + quote-prefix for text cells,
# border for the numeric Quantity column) across the whole table body.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("A2:E17").PasteSpecial(-4122) | Out-Null
$ws.Range("H2").Copy() | Out-Null
$ws.Range("F2:F17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("H1:H2").Clear() | Out-Null
